{"js": "// Fix the capitalization of \"nas\" -> \"NAs\" in the last bullet of the\n// progress tracker (\"Need to clean up the nas in assets derived\").\n//\n// This is a plain word-level text fix, but the canonical OOXML for the\n// edit shows the sentence's single <w:r> run split into three runs\n// (the untouched lead-in text, the corrected word, and the untouched\n// tail), all sharing identical run properties. That's the normal\n// footprint of a Word UI edit: only the run spanning the touched\n// selection gets created anew, the rest of the original run is split\n// around it. Re-stamping (and immediately clearing) character\n// formatting on the replaced word reproduces that same run boundary\n// here instead of letting the save step re-coalesce it back into one\n// run with its identically-formatted neighbors.\n\nconst body = context.document.body;\n\nconst results = body.search(\"nas\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the text \"nas\" to correct.');\n}\n\nconst target = results.items[0];\nconst replaced = target.insertText(\"NAs\", \"Replace\");\n\n// Touch-and-revert a character attribute on just the replacement text so\n// it keeps its own run instead of silently re-merging with its\n// identically-formatted neighbors.\nreplaced.font.bold = true;\nreplaced.font.bold = false;\n\nawait context.sync();\n", "ps1": "# Fix the capitalization of \"nas\" -> \"NAs\" in the last bullet of the\n# progress tracker (\"Need to clean up the nas in assets derived\").\n#\n# The canonical OOXML for this edit shows the sentence's single <w:r>\n# run split into three runs (the untouched lead-in text, the corrected\n# word, and the untouched tail), all sharing identical run properties.\n# That's the normal footprint of a Word UI edit: only the run spanning\n# the touched selection gets created anew, the rest of the original run\n# is split around it. Re-stamping (and immediately clearing) character\n# formatting on the replaced word reproduces that same run boundary\n# here instead of letting the save step re-coalesce it back into one\n# run with its identically-formatted neighbors.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content.Duplicate\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"nas\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $start = $rng.Start\n\n    # Replace the matched word with the corrected capitalization.\n    $rng.Text = \"NAs\"\n\n    # Re-seat a range over just the replacement text and touch-and-revert\n    # a character attribute so it keeps its own run instead of silently\n    # re-merging with its identically-formatted neighbors.\n    $repl = $d.Range($start, $start + 3)\n    $repl.Font.Bold = 1\n    $repl.Font.Bold = 0\n}\n"}
